# Registar planificaciones de la plantilla de planificacion multiple en plan general
#
# The template sheet "Planificacion" held leftover sample/demo data in the
# planning grid (row 3 and rows 4-12). Those cells are cleared out so the
# template starts blank (ready for the real data to be written in by the
# app), while keeping the header rows (1-2) and the cell styles intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planificacion")

# Row 3: keep A3 (day label) and the styled-but-empty B3:D3/E3/F3/J3 cells,
# clear their leftover sample values.
$ws.Range("B3").ClearContents() | Out-Null
$ws.Range("C3").ClearContents() | Out-Null
$ws.Range("D3").ClearContents() | Out-Null
$ws.Range("E3").ClearContents() | Out-Null
$ws.Range("F3").ClearContents() | Out-Null
$ws.Range("J3").ClearContents() | Out-Null

# Rows 4-12: clear the leftover sample values in columns B/C, keep column A.
$ws.Range("B4").ClearContents() | Out-Null

$ws.Range("B5").ClearContents() | Out-Null
$ws.Range("C5").ClearContents() | Out-Null

$ws.Range("C6").ClearContents() | Out-Null
$ws.Range("C7").ClearContents() | Out-Null
$ws.Range("C8").ClearContents() | Out-Null
$ws.Range("C9").ClearContents() | Out-Null
$ws.Range("C10").ClearContents() | Out-Null
$ws.Range("C11").ClearContents() | Out-Null
$ws.Range("C12").ClearContents() | Out-Null

# Move the active selection, matching the author's final cursor position.
$ws.Range("E29").Select() | Out-Null
